$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Estado de Cuenta" data rows (16 and 17):
# - Row 16 ("Periodo Mora") changes from 2109 -> 2108
# - Row 17 ("Periodo Mora") changes from 2108 -> 2109
# - "Valor Mora" values are swapped between the two rows
# - "Salario Basico" is updated to the new value for both rows
$ws.Range("E16").Value = "2108"
$ws.Range("F16").Value = 36341
$ws.Range("G16").Value = 877803

$ws.Range("E17").Value = "2109"
$ws.Range("F17").Value = 35112
$ws.Range("G17").Value = 877803
